$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2100
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2125
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2125
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2475

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4062.5
$ws.Range("I86").Value = 4833.3335
$ws.Range("J86").Value = 3291.6667
$ws.Range("K86").Value = 4833.3335
$ws.Range("L86").Value = 3291.6667
$ws.Range("M86").Value = -3710.3335
$ws.Range("N86").Value = -5537.6667

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4062.5
$ws.Range("I89").Value = 4833.3335
$ws.Range("J89").Value = 3291.6667
$ws.Range("K89").Value = 24166.6675
$ws.Range("L89").Value = 16458.3335
$ws.Range("M89").Value = -18550.6675
$ws.Range("N89").Value = -27690.3335

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 993.625
$ws.Range("I2").Value = 950.41174
$ws.Range("J2").Value = 1098.5714
$ws.Range("K2").Value = 950.41174
$ws.Range("L2").Value = 1098.5714
$ws.Range("M2").Value = -837.41174
$ws.Range("N2").Value = -1324.5714

# Sheet ARM, row 9
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 30000
$ws.Range("J9").Value = 30000
$ws.Range("L9").Value = 30000
$ws.Range("N9").Value = -30340

# Sheet ARM, row 20
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30540

# Sheet ARM, row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -25976

# Sheet ARM, row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 37053
$ws.Range("J55").Value = 37053
$ws.Range("L55").Value = 37053
$ws.Range("N55").Value = -37683

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4766619
$ws.Range("I63").Value = 33333332
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 33333332
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -33332646
$ws.Range("N63").Value = -6872

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4766619
$ws.Range("I66").Value = 33333332
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 166666660
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -166663228
$ws.Range("N66").Value = -34364

# Sheet ARM, row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 18013.5
$ws.Range("J80").Value = 17832.908
$ws.Range("L80").Value = 17832.908
$ws.Range("N80").Value = -19828.908

# Sheet ARM, row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 18013.5
$ws.Range("J83").Value = 17832.908
$ws.Range("L83").Value = 53498.724
$ws.Range("N83").Value = -63482.724

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 13524.154
$ws.Range("I88").Value = 2666.6667
$ws.Range("J88").Value = 22830.572
$ws.Range("K88").Value = 2666.6667
$ws.Range("L88").Value = 22830.572
$ws.Range("M88").Value = -2260.6667
$ws.Range("N88").Value = -23642.572

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 13524.154
$ws.Range("I91").Value = 2666.6667
$ws.Range("J91").Value = 22830.572
$ws.Range("K91").Value = 2666.6667
$ws.Range("L91").Value = 22830.572
$ws.Range("M91").Value = -1262.6667
$ws.Range("N91").Value = -25638.572

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 993.625
$ws.Range("I116").Value = 950.41174
$ws.Range("J116").Value = 1098.5714
$ws.Range("K116").Value = 950.41174
$ws.Range("L116").Value = 1098.5714
$ws.Range("M116").Value = 1343.58826
$ws.Range("N116").Value = -5686.5714

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 993.625
$ws.Range("I3").Value = 950.41174
$ws.Range("J3").Value = 1098.5714
$ws.Range("K3").Value = 950.41174
$ws.Range("L3").Value = 1098.5714
$ws.Range("M3").Value = -836.41174
$ws.Range("N3").Value = -1326.5714

# Sheet BSM, row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Sheet BSM, row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14174.154
$ws.Range("J82").Value = 15668.818
$ws.Range("L82").Value = 15668.818
$ws.Range("N82").Value = -16434.818

# Sheet BSM, row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 14174.154
$ws.Range("J85").Value = 15668.818
$ws.Range("L85").Value = 15668.818
$ws.Range("N85").Value = -18320.818

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1277.7333
$ws.Range("I86").Value = 1276.1428
$ws.Range("J86").Value = 1300
$ws.Range("K86").Value = 1276.1428
$ws.Range("L86").Value = 1300
$ws.Range("M86").Value = -153.1428000000001
$ws.Range("N86").Value = -3546

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1277.7333
$ws.Range("I89").Value = 1276.1428
$ws.Range("J89").Value = 1300
$ws.Range("K89").Value = 6380.714
$ws.Range("L89").Value = 6500
$ws.Range("M89").Value = -764.7139999999999
$ws.Range("N89").Value = -17732

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10575
$ws.Range("I31").Value = 1080.5238
$ws.Range("J31").Value = 23036.5
$ws.Range("K31").Value = 1080.5238
$ws.Range("L31").Value = 23036.5
$ws.Range("M31").Value = -785.5237999999999
$ws.Range("N31").Value = -23626.5

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10575
$ws.Range("I34").Value = 1080.5238
$ws.Range("J34").Value = 23036.5
$ws.Range("K34").Value = 1080.5238
$ws.Range("L34").Value = 23036.5
$ws.Range("M34").Value = -878.5237999999999
$ws.Range("N34").Value = -23440.5

# Sheet CRP, row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Sheet CRP, row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1852.44
$ws.Range("I86").Value = 1536.75
$ws.Range("J86").Value = 2413.6667
$ws.Range("K86").Value = 1536.75
$ws.Range("L86").Value = 2413.6667
$ws.Range("M86").Value = -413.75
$ws.Range("N86").Value = -4659.6667

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 1852.44
$ws.Range("I89").Value = 1536.75
$ws.Range("J89").Value = 2413.6667
$ws.Range("K89").Value = 7683.75
$ws.Range("L89").Value = 12068.3335
$ws.Range("M89").Value = -2067.75
$ws.Range("N89").Value = -23300.3335

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 348.0625
$ws.Range("I107").Value = 208.57143
$ws.Range("J107").Value = 456.55554
$ws.Range("K107").Value = 208.57143
$ws.Range("L107").Value = 456.55554
$ws.Range("M107").Value = 1711.42857
$ws.Range("N107").Value = -4296.55554

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 720.8
$ws.Range("I5").Value = 651
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1953
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1841
$ws.Range("N5").Value = -3224

# Sheet CUL, row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 923.4286
$ws.Range("J117").Value = 923.4286
$ws.Range("L117").Value = 2770.2858
$ws.Range("N117").Value = -9654.2858

# Sheet CUL, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 822.36
$ws.Range("I129").Value = 351
$ws.Range("J129").Value = 1087.5
$ws.Range("K129").Value = 1053
$ws.Range("L129").Value = 3262.5
$ws.Range("M129").Value = 3947
$ws.Range("N129").Value = -13262.5

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 720.8
$ws.Range("I135").Value = 651
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 5859
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3324
$ws.Range("N135").Value = -14070

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2585.5
$ws.Range("I80").Value = 2425.8333
$ws.Range("J80").Value = 2825
$ws.Range("K80").Value = 2425.8333
$ws.Range("L80").Value = 2825
$ws.Range("M80").Value = -1427.8333
$ws.Range("N80").Value = -4821

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2585.5
$ws.Range("I83").Value = 2425.8333
$ws.Range("J83").Value = 2825
$ws.Range("K83").Value = 12129.1665
$ws.Range("L83").Value = 14125
$ws.Range("M83").Value = -7137.166499999999
$ws.Range("N83").Value = -24109

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 50600
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 50600
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 50600
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -50940

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5541
$ws.Range("I46").Value = 1083.4445
$ws.Range("J46").Value = 25600
$ws.Range("K46").Value = 1083.4445
$ws.Range("L46").Value = 25600
$ws.Range("M46").Value = -895.4445000000001
$ws.Range("N46").Value = -25976

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20836258
$ws.Range("I81").Value = 1462
$ws.Range("K81").Value = 2924
$ws.Range("M81").Value = -1863

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 20836258
$ws.Range("I84").Value = 1462
$ws.Range("K84").Value = 14620
$ws.Range("M84").Value = -9316

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21403230
$ws.Range("I136").Value = 11128412
$ws.Range("J136").Value = 62502500
$ws.Range("K136").Value = 33385236
$ws.Range("L136").Value = 187507500
$ws.Range("M136").Value = -33382686
$ws.Range("N136").Value = -187512600

